$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3196676666666667
$ws.Range("H2").Value = 0.959003
$ws.Range("I2").Value = 0.02778181741070332
$ws.Range("J2").Value = 0.02778181741070332
$ws.Range("M2").Value = 110.642708
$ws.Range("N2").Value = 331.928124
$ws.Range("O2").Value = 0.5476418925386564
$ws.Range("P2").Value = 0.5476418925386564
$ws.Range("Q2").Value = 35.36889630004133
$ws.Range("R2").Value = 318.320066700372
$ws.Range("S2").Value = 0.01521448706496096
$ws.Range("T2").Value = 0.01521448706496096
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3196676666666667
$ws.Range("H3").Value = 0.959003
$ws.Range("I3").Value = 0.02778181741070332
$ws.Range("J3").Value = 0.02778181741070332
$ws.Range("O3").Value = 0.3151072754333865
$ws.Range("P3").Value = 0.3151072754333865
$ws.Range("Q3").Value = 20.35088385318389
$ws.Range("R3").Value = 183.157954678655
$ws.Range("S3").Value = 0.008754252790874543
$ws.Range("T3").Value = 0.008754252790874545
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3196676666666667
$ws.Range("H4").Value = 0.959003
$ws.Range("I4").Value = 0.02778181741070332
$ws.Range("J4").Value = 0.02778181741070332
$ws.Range("M4").Value = 27.72944133333333
$ws.Range("N4").Value = 83.18832399999999
$ws.Range("O4").Value = 0.1372508320279571
$ws.Range("P4").Value = 0.1372508320279571
$ws.Range("Q4").Value = 8.864205808996889
$ws.Range("R4").Value = 79.777852280972
$ws.Range("S4").Value = 0.003813077554867816
$ws.Range("T4").Value = 0.003813077554867816
$ws.Range("G5").Value = 5.787456
$ws.Range("I5").Value = 0.5029787577238425
$ws.Range("J5").Value = 0.5029787577238426
$ws.Range("M5").Value = 110.642708
$ws.Range("N5").Value = 331.928124
$ws.Range("O5").Value = 0.5476418925386564
$ws.Range("P5").Value = 0.5476418925386564
$ws.Range("Q5").Value = 640.3398042708479
$ws.Range("R5").Value = 5763.058238437631
$ws.Range("S5").Value = 0.2754522387866274
$ws.Range("T5").Value = 0.2754522387866275
$ws.Range("G6").Value = 5.787456
$ws.Range("I6").Value = 0.5029787577238425
$ws.Range("J6").Value = 0.5029787577238426
$ws.Range("O6").Value = 0.3151072754333865
$ws.Range("P6").Value = 0.3151072754333865
$ws.Range("Q6").Value = 368.44466032352
$ws.Range("S6").Value = 0.1584922659472294
$ws.Range("T6").Value = 0.1584922659472295
$ws.Range("G7").Value = 5.787456
$ws.Range("I7").Value = 0.5029787577238425
$ws.Range("J7").Value = 0.5029787577238426
$ws.Range("M7").Value = 27.72944133333333
$ws.Range("N7").Value = 83.18832399999999
$ws.Range("O7").Value = 0.1372508320279571
$ws.Range("P7").Value = 0.1372508320279571
$ws.Range("Q7").Value = 160.482921621248
$ws.Range("R7").Value = 1444.346294591232
$ws.Range("S7").Value = 0.06903425298998565
$ws.Range("T7").Value = 0.06903425298998567
$ws.Range("G8").Value = 5.399239000000001
$ws.Range("H8").Value = 16.197717
$ws.Range("I8").Value = 0.4692394248654542
$ws.Range("J8").Value = 0.4692394248654542
$ws.Range("M8").Value = 110.642708
$ws.Range("N8").Value = 331.928124
$ws.Range("O8").Value = 0.5476418925386564
$ws.Range("P8").Value = 0.5476418925386564
$ws.Range("Q8").Value = 597.386424099212
$ws.Range("R8").Value = 5376.477816892908
$ws.Range("S8").Value = 0.256975166687068
$ws.Range("T8").Value = 0.256975166687068
$ws.Range("G9").Value = 5.399239000000001
$ws.Range("H9").Value = 16.197717
$ws.Range("I9").Value = 0.4692394248654542
$ws.Range("J9").Value = 0.4692394248654542
$ws.Range("O9").Value = 0.3151072754333865
$ws.Range("P9").Value = 0.3151072754333865
$ws.Range("Q9").Value = 343.7297457398384
$ws.Range("R9").Value = 3093.567711658545
$ws.Range("S9").Value = 0.1478607566952825
$ws.Range("T9").Value = 0.1478607566952826
$ws.Range("G10").Value = 5.399239000000001
$ws.Range("H10").Value = 16.197717
$ws.Range("I10").Value = 0.4692394248654542
$ws.Range("J10").Value = 0.4692394248654542
$ws.Range("M10").Value = 27.72944133333333
$ws.Range("N10").Value = 83.18832399999999
$ws.Range("O10").Value = 0.1372508320279571
$ws.Range("P10").Value = 0.1372508320279571
$ws.Range("Q10").Value = 149.7178810951453
$ws.Range("R10").Value = 1347.460929856308
$ws.Range("S10").Value = 0.06440350148310366
$ws.Range("T10").Value = 0.06440350148310367
